$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 20 data: a new task entry "Работа по созданию функционала загрузки данных (Добавление офтальмологии)"
$ws.Range("A20").Value = "Работа по созданию функционала загрузки данных (Добавление офтальмологии)"
$ws.Range("B20").Value = 1
$ws.Range("C19").Copy($ws.Range("C20"))
$ws.Range("C20").Value = 43552

# Update the selection to match the new active cell
$ws.Range("C21").Select()
